# DataArray und TestDataArray eingefügt
# ---------------------------------------------------------------
# Semantic edits applied (per the commit):
#  - Mark "Objekte als Rückgabewerte" (B11/row11) and
#    "Kombinationsaufgabe aus allen Themen" (B12/row12) as done
#    (checkbox glyph "a" in the Marlett font, same as the other
#    checked rows in the list).
#  - Document, in the "Bemerkungen" column (E), that row 11's
#    topic was already covered together with Data3, and that
#    row 12 covers "DataArray".
#  - Re-balance the light-grey row striping for the "Data3" /
#    "Shallow-Copy & Deep-Copy" rows (13/14) and draw a small
#    left-hand bracket connecting them, since they visually
#    belong together now that row 11 references Data3.
#  - Row 11 grows (its note now wraps), selection cursor ends on H13.
# ---------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tick the checkboxes on rows 11 and 12 ------------------------------
$ws.Range("B11").Value = "a"
$ws.Range("B12").Value = "a"

# --- New remarks in column E --------------------------------------------
$ws.Range("E11").Value = "wurde in Data3 mit behandelt"
$ws.Range("E12").Value = "DataArray"

# --- Row 11 grows a bit now that it carries a longer remark -------------
$ws.Rows.Item(11).RowHeight = 25

# --- Re-align the alternating row shading for rows 13/14 -----------------
# Row 13 ("Data3") takes over the "shaded" look that rows 7/9/11 use,
# row 14 ("Shallow-Copy & Deep-Copy") takes over the "plain" look that
# rows 25/36 use - this keeps the stripes alternating correctly now that
# row 11 grew.
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null

$ws.Range("C8").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null

$ws.Range("D9").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null

$ws.Range("D25").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null

$ws.Range("E9").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null

$ws.Range("E25").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Small connecting bracket between rows 13 and 14 on column C --------
$ws.Range("C13").Borders.Item(7).LineStyle = 1
$ws.Range("C14").Borders.Item(7).LineStyle = 1
$ws.Range("C14").Borders.Item(9).LineStyle = 1

# --- Cursor position, matching where editing finished --------------------
$ws.Range("H13").Select() | Out-Null
